$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial number. Every data row (2-416)
# had its value bumped by one day: 45203 (2023-10-04) -> 45204 (2023-10-05).
$lastRow = 416

for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    $current = $cell.Value2()
    if ($current -eq 45203) {
        $cell.Value = 45204
    }
}
